$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 0.8229573333333334
$ws.Range("N2").Value = 2.468872
$ws.Range("O2").Value = 0.2440777672676426
$ws.Range("P2").Value = 0.2440777672676426
$ws.Range("Q2").Value = 0.3735200339857778
$ws.Range("R2").Value = 3.361680305872
$ws.Range("S2").Value = 0.2440777672676426
$ws.Range("T2").Value = 0.2440777672676426

$ws.Range("O3").Value = 0.4345811965947162
$ws.Range("P3").Value = 0.4345811965947162
$ws.Range("S3").Value = 0.4345811965947162
$ws.Range("T3").Value = 0.4345811965947162

$ws.Range("M4").Value = 1.083466
$ws.Range("N4").Value = 3.250398
$ws.Range("O4").Value = 0.3213410361376413
$ws.Range("P4").Value = 0.3213410361376413
$ws.Range("Q4").Value = 0.4917584919053333
$ws.Range("R4").Value = 4.425826427147999
$ws.Range("S4").Value = 0.3213410361376413
$ws.Range("T4").Value = 0.3213410361376413
